$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they stay as text (matching source data),
# rather than being auto-converted to numbers by Excel. (Looping since multi-area Range
# assignment only reliably touches the first area in this host.)
$textCells = @("D5", "D6", "D14", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D29", "D30", "D35", "D36", "D39", "D40", "D42", "D43", "D46", "D50")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Updated Price / Volume(1h) figures
$ws.Range("D2").Value = "60.882.48"
$ws.Range("E2").Value = "  +5.64%  "
$ws.Range("D3").Value = "2.370.67"
$ws.Range("E3").Value = "  +4.15%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "548.96"
$ws.Range("E5").Value = "  +3.64%  "
$ws.Range("D6").Value = "133.10"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "2.368.44"
$ws.Range("E9").Value = "  +3.91%  "
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "24.13"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "2.795.42"
$ws.Range("E15").Value = "  +4.25%  "
$ws.Range("D16").Value = "60.820.12"
$ws.Range("E16").Value = "  +5.51%  "
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "2.375.75"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("D19").Value = "10.77"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").Value = "6.96"
$ws.Range("E20").Value = "  +9.99%  "
$ws.Range("D21").Value = "4.20"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "317.82"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "63.56"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("E25").Value = "  +4.24%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "8.05"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").Value = "1.75"
$ws.Range("E29").Value = "  +3.38%  "
$ws.Range("D30").Value = "172.16"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "0.0₃0738"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("E32").Value = "  +11.59%  "
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("E34").Value = "  +16.82%  "
$ws.Range("D35").Value = "0.385"
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").Value = "18.12"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "4.20"
$ws.Range("E39").Value = "  +8.76%  "
$ws.Range("D40").Value = "319.53"
$ws.Range("E40").Value = "  +12.00%  "
$ws.Range("E41").Value = "  +4.67%  "
$ws.Range("D42").Value = "38.29"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "144.58"
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "19.39"
$ws.Range("E46").Value = "  +8.03%  "
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("E49").Value = "  +2.90%  "

# Rows 50 and 51 swapped order (WhiteBITCoin now ranked above BabyDogeCoin), with updated values
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "11.04"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0209"
$ws.Range("E51").Value = "  +4.85%  "
